# Add the missing week-15 matchup (Boomer Sooners vs. Kuppenheimer) that was
# left out of the ladder when there was a bye, by inserting two new rows right
# before the week-16 matchups (row 182) and filling them with the matchup data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the week-16+ rows down by two to make room for the new matchup.
$ws.Rows("182:183").Insert()

# Row 182: Boomer Sooners (Home) vs. Kuppenheimer
$ws.Cells.Item(182, 1).Value = "Boomer Sooners"
$ws.Cells.Item(182, 2).Value = "Kuppenheimer"
$ws.Cells.Item(182, 3).Value = 15
$ws.Cells.Item(182, 4).Value = "Home"
$ws.Cells.Item(182, 5).Value = 92.08000000000001
$ws.Cells.Item(182, 6).Value = 132.44
$ws.Cells.Item(182, 7).Value = "LOSERS_CONSOLATION_LADDER"
$ws.Cells.Item(182, 8).Value = -1.99
$ws.Cells.Item(182, 9).Value = 63.28
$ws.Cells.Item(182, 10).Value = "Boomer Sooners vs. Kuppenheimer"
$ws.Cells.Item(182, 11).Value = "Win"
$ws.Cells.Item(182, 12).Value = 2
$ws.Cells.Item(182, 13).Value = "15Boomer Sooners"
$ws.Cells.Item(182, 14).Value = 7
$ws.Cells.Item(182, 15).Value = 13

# Row 183: Kuppenheimer (Away) vs. Boomer Sooners
$ws.Cells.Item(183, 1).Value = "Kuppenheimer"
$ws.Cells.Item(183, 2).Value = "Boomer Sooners"
$ws.Cells.Item(183, 3).Value = 15
$ws.Cells.Item(183, 4).Value = "Away"
$ws.Cells.Item(183, 5).Value = 94.07000000000001
$ws.Cells.Item(183, 6).Value = 69.16
$ws.Cells.Item(183, 7).Value = "LOSERS_CONSOLATION_LADDER"
$ws.Cells.Item(183, 8).Value = 1.99
$ws.Cells.Item(183, 9).Value = -63.28
$ws.Cells.Item(183, 10).Value = "Boomer Sooners vs. Kuppenheimer"
$ws.Cells.Item(183, 11).Value = "Loss"
$ws.Cells.Item(183, 12).Value = 1
$ws.Cells.Item(183, 13).Value = "15Kuppenheimer"
$ws.Cells.Item(183, 14).Value = 7
$ws.Cells.Item(183, 15).Value = 7
